$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.877.22"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.450.83"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.83"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.00"
$ws.Range("E6").Value = "  +7.05%  "
$ws.Range("D7").Value = "3.451.41"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.64"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "4.043.46"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.67"
$ws.Range("E15").Value = "  +6.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "3.455.40"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "62.010.55"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.30"
$ws.Range("E19").Value = "  +6.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.14"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.53"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.57"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.563"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "3.593.23"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.78"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "72.43"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000124"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  +8.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.77"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.54"
$ws.Range("E31").Value = "  -15.37%  "
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.25"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.98"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.28"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.04"
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.19"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0791"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.00"
$ws.Range("E42").Value = "  +7.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.793"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.73"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.12"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.661.02"
$ws.Range("E48").Value = "  +11.92%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.91"
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.86"
$ws.Range("E51").Value = "  -0.01%  "
